# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E).
# Values are plain text in the source sheet (t="inlineStr"), some of which read as
# plain numbers (e.g. "0.9998"); force NumberFormat "@" first so Excel keeps them as
# text instead of silently coercing them to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.212.02'; E = '  +0.70%  ' },
    @{ Row = 3; D = '1.901.43'; E = '  +0.65%  ' },
    @{ Row = 4; D = $null; E = '  -0.08%  ' },
    @{ Row = 5; D = '307.76'; E = '  +0.53%  ' },
    @{ Row = 6; D = '0.9998'; E = '  -0.08%  ' },
    @{ Row = 7; D = '0.5205'; E = '  +0.30%  ' },
    @{ Row = 8; D = '0.3772'; E = '  +0.52%  ' },
    @{ Row = 9; D = '0.07275'; E = '  +1.06%  ' },
    @{ Row = 10; D = '21.18'; E = '  +0.26%  ' },
    @{ Row = 11; D = '0.9033'; E = '  +0.23%  ' },
    @{ Row = 12; D = '0.08272'; E = '  +8.39%  ' },
    @{ Row = 13; D = '1.917.31'; E = '  +2.12%  ' },
    @{ Row = 14; D = '96.43'; E = '  +2.07%  ' },
    @{ Row = 15; D = '5.277'; E = '  +0.73%  ' },
    @{ Row = 16; D = '1.000'; E = '  -0.14%  ' },
    @{ Row = 17; D = '0.000008636'; E = '  +1.58%  ' },
    @{ Row = 18; D = '14.56'; E = '  +0.78%  ' },
    @{ Row = 19; D = '0.9998'; E = '  -0.02%  ' },
    @{ Row = 20; D = '27.233.68'; E = '  +0.61%  ' },
    @{ Row = 21; D = '5.089'; E = '  +0.76%  ' },
    @{ Row = 22; D = '2.158.30'; E = '  +1.70%  ' },
    @{ Row = 23; D = $null; E = '  +0.62%  ' },
    @{ Row = 24; D = '6.426'; E = '  +0.60%  ' },
    @{ Row = 25; D = '2.323'; E = '  +1.11%  ' },
    @{ Row = 26; D = '147.11'; E = '  +0.90%  ' },
    @{ Row = 27; D = '1.747'; E = '  +0.65%  ' },
    @{ Row = 28; D = $null; E = '  +0.80%  ' },
    @{ Row = 29; D = '115.01'; E = '  +0.73%  ' },
    @{ Row = 30; D = '4.840'; E = '  +1.14%  ' },
    @{ Row = 31; D = '4.902'; E = '  -0.26%  ' },
    @{ Row = 32; D = '0.09247'; E = '  +0.61%  ' },
    @{ Row = 33; D = '0.05074'; E = '  +0.66%  ' },
    @{ Row = 34; D = '0.7977'; E = '  +4.03%  ' },
    @{ Row = 35; D = '1.237'; E = '  -0.18%  ' },
    @{ Row = 36; D = $null; E = '  +4.61%  ' },
    @{ Row = 37; D = '2.956'; E = '  -0.16%  ' },
    @{ Row = 38; D = '2.599'; E = '  -0.19%  ' },
    @{ Row = 39; D = '0.5709'; E = '  +1.98%  ' },
    @{ Row = 40; D = '0.02002'; E = '  +0.72%  ' },
    @{ Row = 41; D = $null; E = '  +0.38%  ' },
    @{ Row = 42; D = '9.034'; E = $null },
    @{ Row = 43; D = '6.592'; E = '  -0.37%  ' },
    @{ Row = 44; D = '116.70'; E = '  -1.74%  ' },
    @{ Row = 45; D = '0.1517'; E = '  +0.80%  ' },
    @{ Row = 46; D = '0.4857'; E = '  +0.67%  ' },
    @{ Row = 47; D = '0.9998'; E = '  -0.05%  ' },
    @{ Row = 48; D = '10.07'; E = '  -1.15%  ' },
    @{ Row = 49; D = '1.629'; E = '  +2.08%  ' },
    @{ Row = 50; D = '37.64'; E = '  -0.11%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$($u.Row)")
        if ($u.D -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
